# Updated cryptos list on Sat Nov 11 16:25:04 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = '37.152.53'
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").Value = '2.072.42'
$ws.Range("E3").Value = '  -1.28%  '
$ws.Range("E4").Value = '  +0.11%  '
Set-TextValue "D5" '253.97'
$ws.Range("E5").Value = '  +1.43%  '
$ws.Range("E6").Value = '  +3.18%  '
Set-TextValue "D7" '61.83'
$ws.Range("E7").Value = '  +20.65%  '
Set-TextValue "D8" '1.00'
$ws.Range("E8").Value = '  +0.00%  '
Set-TextValue "D9" '0.394'
$ws.Range("E9").Value = '  +5.97%  '
Set-TextValue "D10" '61.63'
$ws.Range("E10").Value = '  -0.75%  '
Set-TextValue "D11" '0.0810'
$ws.Range("E11").Value = '  +9.12%  '
Set-TextValue "D12" '0.109'
$ws.Range("E12").Value = '  +3.12%  '
Set-TextValue "D13" '16.46'
$ws.Range("E13").Value = '  +7.25%  '
$ws.Range("D14").Value = '2.375.54'
$ws.Range("E14").Value = '  -1.16%  '
Set-TextValue "D15" '0.823'
$ws.Range("E15").Value = '  -0.68%  '
Set-TextValue "D16" '5.54'
$ws.Range("D17").Value = '2.070.97'
$ws.Range("E17").Value = '  -1.19%  '
$ws.Range("D18").Value = '37.099.07'
$ws.Range("E18").Value = '  -0.28%  '
Set-TextValue "D19" '74.85'
$ws.Range("E19").Value = '  +3.84%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D20" '15.50'
$ws.Range("E20").Value = '  +14.56%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.0₃0932'
$ws.Range("E21").Value = '  +12.69%  '
Set-TextValue "D22" '5.51'
$ws.Range("E22").Value = '  +5.57%  '
Set-TextValue "D23" '240.36'
$ws.Range("E23").Value = '  +0.22%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").Value = '  -0.95%  '
Set-TextValue "D26" '2.34'
$ws.Range("E26").Value = '  +17.64%  '
Set-TextValue "D27" '170.42'
$ws.Range("E27").Value = '  +0.41%  '
Set-TextValue "D28" '9.34'
$ws.Range("E28").Value = '  +2.05%  '
Set-TextValue "D29" '20.43'
$ws.Range("E29").Value = '  -0.76%  '
Set-TextValue "D30" '0.127'
$ws.Range("E30").Value = '  +3.62%  '
Set-TextValue "D31" '4.86'
$ws.Range("E31").Value = '  +8.32%  '
$ws.Range("E32").Value = '  +5.09%  '
Set-TextValue "D33" '0.0640'
$ws.Range("E33").Value = '  +5.49%  '
Set-TextValue "D34" '4.46'
$ws.Range("E34").Value = '  +9.58%  '
Set-TextValue "D35" '0.0898'
$ws.Range("E35").Value = '  -2.49%  '
Set-TextValue "D36" '1.00'
$ws.Range("E36").Value = '  +0.06%  '
Set-TextValue "D37" '2.31'
$ws.Range("E37").Value = '  -0.87%  '
$ws.Range("E38").Value = '  -2.92%  '
Set-TextValue "D39" '0.113'
$ws.Range("E39").Value = '  +25.98%  '
Set-TextValue "D40" '1.38'
$ws.Range("E40").Value = '  +5.42%  '
Set-TextValue "D41" '18.05'
$ws.Range("E41").Value = '  +1.21%  '
Set-TextValue "D42" '0.0228'
$ws.Range("E42").Value = '  +2.01%  '
Set-TextValue "D43" '1.17'
$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue "D44" '4.49'
$ws.Range("E44").Value = '  +28.65%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D45" '98.99'
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("E46").Value = '  +2.94%  '
Set-TextValue "D47" '4.57'
$ws.Range("E47").Value = '  +13.67%  '
$ws.Range("E48").Value = '  +10.28%  '
$ws.Range("D49").Value = '1.306.08'
$ws.Range("E49").Value = '  -1.00%  '
$ws.Range("E50").Value = '  -2.36%  '
Set-TextValue "D51" '6.94'
$ws.Range("E51").Value = '  -0.12%  '
